# "mejoras de clientes: dialog, plantillas, emailjs"
# Adds a new "logoCliente" column (K) to the Usuarios sheet with a
# logo-image URL per client, each one a hyperlink, matching the styling
# of the existing "sitioWeb"/"correo" hyperlink columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header + data text, entered in this order so the shared-string
# table is built up with the same indices as the authored workbook.
$ws.Cells.Item(1, 11).Value = "logoCliente"
$ws.Cells.Item(2, 11).Value = "https://www.ivelpink.cl/logo-ivelpink-correo.png"
$ws.Cells.Item(4, 11).Value = "https://sifg.cl/logo-sifg-correo.png"
$ws.Cells.Item(3, 11).Value = "https://masautomatizacion.cl/logo-masautomatizacion.jpg"
$ws.Cells.Item(5, 11).Value = "https://ingsnt.cl/logo-ingsnt-correo-white.webp"

# --- Turn each logo URL into a hyperlink. Added in this order
# (masautomatizacion, ivelpink, sifg, ingsnt) to reproduce the
# relationship-id allocation order (rId13..rId16).
$ws.Hyperlinks.Add($ws.Range("K3"), "https://masautomatizacion.cl/logo-masautomatizacion.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K2"), "https://www.ivelpink.cl/logo-ivelpink-correo.png") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K4"), "https://sifg.cl/logo-sifg-correo.png") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K5"), "https://ingsnt.cl/logo-ingsnt-correo-white.webp") | Out-Null

# --- Formatting -----------------------------------------------------
# Header cell: match the other header cells (bold / centered / wrap),
# i.e. style index 1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Hyperlink cells: match the formatting already used for the other
# hyperlink cells (e.g. C2, which uses the shared "Hyperlink" cell
# style) instead of the new style Hyperlinks.Add() would synthesize.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Column width for the new column K --------------------------------
$ws.Columns.Item(11).ColumnWidth = 49.8

# --- Selection shown when the file was last saved ----------------------
$ws.Range("K6").Select()
